# Append two new tracker rows (01.04.2024 report + 02.04.2024 tracker update)
# to the mortality tracker sheet, mirroring the pattern used for prior rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 177
$row2 = 178

function Set-TextValue($cell, $text) {
    # Force the cell to keep its value as literal text (not auto-parsed
    # into a date serial), then restore the default "General" number
    # format so the cell style matches the rest of the row.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.NumberFormat = "General"
}

# --- Row 177: new report_date 01.04.2024, tracker_date also 01.04.2024 ---
Set-TextValue $ws.Cells.Item($row1, 1) "01.04.2024"
Set-TextValue $ws.Cells.Item($row1, 2) "01.04.2024"
$ws.Cells.Item($row1, 3).Value = 32845
$ws.Cells.Item($row1, 4).Value = 13000
$ws.Cells.Item($row1, 5).Value = 8400
$ws.Cells.Item($row1, 6).Value = 75392
$ws.Cells.Item($row1, 7).Value = 8663
$ws.Cells.Item($row1, 8).Value = 6327
$ws.Cells.Item($row1, 9).Value = 8000
$ws.Cells.Item($row1, 10).Value = 456
$ws.Cells.Item($row1, 11).Value = 117
$ws.Cells.Item($row1, 12).Value = 4750
Set-TextValue $ws.Cells.Item($row1, 13) "https://web.archive.org/web/20240402082823/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 178: new tracker_date 02.04.2024, report_date still 01.04.2024 ---
Set-TextValue $ws.Cells.Item($row2, 1) "02.04.2024"
Set-TextValue $ws.Cells.Item($row2, 2) "01.04.2024"
$ws.Cells.Item($row2, 3).Value = 32845
$ws.Cells.Item($row2, 4).Value = 13000
$ws.Cells.Item($row2, 5).Value = 8400
$ws.Cells.Item($row2, 6).Value = 75392
$ws.Cells.Item($row2, 7).Value = 8663
$ws.Cells.Item($row2, 8).Value = 6327
$ws.Cells.Item($row2, 9).Value = 8000
$ws.Cells.Item($row2, 10).Value = 456
$ws.Cells.Item($row2, 11).Value = 117
$ws.Cells.Item($row2, 12).Value = 4750
Set-TextValue $ws.Cells.Item($row2, 13) "https://web.archive.org/web/20240402173216/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"
